# Apply crypto price/volume updates (cryptos list refresh, GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "98.034.27"
Set-TextCell 2 5 "  -0.36%  "

Set-TextCell 3 4 "3.379.92"
Set-TextCell 3 5 "  -1.15%  "

Set-TextCell 4 4 "0.999"
Set-TextCell 4 5 "  -0.07%  "

Set-TextCell 5 4 "252.81"
Set-TextCell 5 5 "  -1.39%  "

Set-TextCell 6 4 "658.96"
Set-TextCell 6 5 "  -1.43%  "

Set-TextCell 7 4 "1.48"
Set-TextCell 7 5 "  +0.94%  "

Set-TextCell 8 4 "0.427"
Set-TextCell 8 5 "  -2.31%  "

Set-TextCell 9 5 "  -0.06%  "

Set-TextCell 10 4 "1.04"
Set-TextCell 10 5 "  -2.90%  "

Set-TextCell 11 4 "3.376.18"
Set-TextCell 11 5 "  -1.19%  "

Set-TextCell 12 5 "  -3.18%  "

Set-TextCell 13 4 "43.56"
Set-TextCell 13 5 "  +3.16%  "

Set-TextCell 14 4 "97.754.35"
Set-TextCell 14 5 "  -0.36%  "

Set-TextCell 15 4 "6.11"
Set-TextCell 15 5 "  -5.29%  "

Set-TextCell 16 4 "0.0000258"
Set-TextCell 16 5 "  -3.80%  "

Set-TextCell 17 4 "4.008.32"
Set-TextCell 17 5 "  -1.06%  "

Set-TextCell 18 4 "9.24"
Set-TextCell 18 5 "  +1.76%  "

Set-TextCell 19 4 "3.384.61"
Set-TextCell 19 5 "  -0.57%  "

Set-TextCell 20 4 "18.05"
Set-TextCell 20 5 "  +1.85%  "

Set-TextCell 21 4 "0.524"
Set-TextCell 21 5 "  -11.59%  "

Set-TextCell 22 4 "11.36"
Set-TextCell 22 5 "  +2.69%  "

Set-TextCell 23 4 "508.95"
Set-TextCell 23 5 "  -0.68%  "

Set-TextCell 24 4 "3.41"
Set-TextCell 24 5 "  -1.42%  "

Set-TextCell 25 5 "  -3.06%  "

Set-TextCell 26 5 "  +3.26%  "

Set-TextCell 27 4 "96.72"
Set-TextCell 27 5 "  -4.99%  "

Set-TextCell 28 4 "12.37"
Set-TextCell 28 5 "  -3.84%  "

Set-TextCell 29 4 "3.558.71"
Set-TextCell 29 5 "  -1.13%  "

Set-TextCell 30 4 "11.69"
Set-TextCell 30 5 "  +0.49%  "

Set-TextCell 31 5 "  -5.81%  "

Set-TextCell 32 5 "  -0.11%  "

Set-TextCell 33 4 "0.191"
Set-TextCell 33 5 "  -3.38%  "

Set-TextCell 34 4 "2.63"
Set-TextCell 34 5 "  +4.95%  "

Set-TextCell 35 4 "0.998"
Set-TextCell 35 5 "  -0.11%  "

Set-TextCell 36 4 "0.562"
Set-TextCell 36 5 "  -2.52%  "

Set-TextCell 37 4 "28.80"
Set-TextCell 37 5 "  -4.35%  "

Set-TextCell 38 5 "  -1.00%  "

Set-TextCell 39 5 "  -2.15%  "

Set-TextCell 40 4 "522.79"
Set-TextCell 40 5 "  -2.99%  "

Set-TextCell 41 5 "  -1.65%  "

Set-TextCell 42 5 "  -0.02%  "

Set-TextCell 43 4 "24.40"
Set-TextCell 43 5 "  -1.30%  "

Set-TextCell 44 2 "ARBITRUM"
Set-TextCell 44 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell 44 4 "0.849"
Set-TextCell 44 5 "  -3.52%  "

Set-TextCell 45 2 "ImmutableX"
Set-TextCell 45 3 "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell 45 4 "1.74"
Set-TextCell 45 5 "  +0.33%  "

Set-TextCell 46 5 "  -2.39%  "

Set-TextCell 47 4 "3.69"
Set-TextCell 47 5 "  -3.38%  "

Set-TextCell 48 2 "Cosmos"
Set-TextCell 48 3 "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell 48 4 "8.67"
Set-TextCell 48 5 "  -3.53%  "

Set-TextCell 49 2 "Stacks"
Set-TextCell 49 3 "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell 49 4 "2.25"
Set-TextCell 49 5 "  +5.83%  "

Set-TextCell 50 4 "5.59"
Set-TextCell 50 5 "  -5.13%  "

Set-TextCell 51 4 "55.46"
Set-TextCell 51 5 "  +2.58%  "
